$d = $word.ActiveDocument
$d.Content.Find.Execute("AMS core component", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PnP core component", 2)
